$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31 -------------------------------------------------------------
# event_id (text, not numeric) and fecha (text, not date) need the leading
# apostrophe so the engine stores them as literal text instead of coercing
# to a number / date serial. ClearFormats() drops the quote-prefix style
# that gets attached so the cell keeps its default (unstyled) look.
$ws.Range("A31").Value = "'14655141"
$ws.Range("A31").ClearFormats()

$ws.Range("B31").Value = "'2025-09-20"
$ws.Range("B31").ClearFormats()

$ws.Range("C31").Value = "Michael Zheng"
$ws.Range("D31").Value = "Daniel Masur"
$ws.Range("E31").Value = "Gana Daniel Masur"

$ws.Range("F31").Value = 2.75

# resultado / profit are still blank (pending match) - write as empty text
# cells (leading apostrophe -> empty string) rather than leaving them
# completely unset, matching the "blank but present" cells used elsewhere
# in the sheet.
$ws.Range("G31").Value = "'"
$ws.Range("G31").ClearFormats()
$ws.Range("H31").Value = "'"
$ws.Range("H31").ClearFormats()

# --- Row 32 -------------------------------------------------------------
$ws.Range("A32").Value = "'14725672"
$ws.Range("A32").ClearFormats()

$ws.Range("B32").Value = "'2025-09-20"
$ws.Range("B32").ClearFormats()

$ws.Range("C32").Value = "Kris van Wyk"
$ws.Range("D32").Value = "Liam Broady"
$ws.Range("E32").Value = "Gana Kris van Wyk"

$ws.Range("F32").Value = 3.25

$ws.Range("G32").Value = "'"
$ws.Range("G32").ClearFormats()
$ws.Range("H32").Value = "'"
$ws.Range("H32").ClearFormats()
